$d = $word.ActiveDocument

$replacements = @(
    @("15×27=405", "56×94=5264"),
    @("13×40=520", "32×62=1984"),
    @("72×25=1800", "46×46=2116"),
    @("31×53=1643", "34×30=1020"),
    @("59×81=4779", "25×44=1100"),
    @("78×66=5148", "31×32=992"),
    @("15×94=1410", "62×59=3658"),
    @("81×89=7209", "12×33=396"),
    @("54×32=1728", "72×38=2736"),
    @("21×95=1995", "51×69=3519"),
    @("95×68=6460", "41×74=3034"),
    @("87×17=1479", "89×15=1335"),
    @("34×80=2720", "22×64=1408"),
    @("98×16=1568", "56×89=4984"),
    @("16×52=832", "60×56=3360"),
    @("54×31=1674", "47×57=2679"),
    @("24×32=768", "29×98=2842"),
    @("26×29=754", "61×35=2135"),
    @("36×74=2664", "24×15=360"),
    @("54×56=3024", "42×58=2436"),
    @("48×17=816", "49×25=1225"),
    @("63×65=4095", "87×30=2610"),
    @("28×12=336", "62×79=4898"),
    @("59×77=4543", "28×44=1232"),
    @("99×47=4653", "34×94=3196")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
